# Update generated "想去人数" (interest count) figures in column F
# for the "展览" and "全部类型" sheets, as refreshed by the gh-pages build.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): rows 2-35 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3127
$ws1.Range("F4").Value = 1088
$ws1.Range("F6").Value = 30
$ws1.Range("F8").Value = 34
$ws1.Range("F9").Value = 1120
$ws1.Range("F10").Value = 15577
$ws1.Range("F11").Value = 234
$ws1.Range("F12").Value = 168
$ws1.Range("F13").Value = 1022
$ws1.Range("F14").Value = 6142
$ws1.Range("F15").Value = 619
$ws1.Range("F22").Value = 630
$ws1.Range("F23").Value = 9
$ws1.Range("F24").Value = 9
$ws1.Range("F26").Value = 207
$ws1.Range("F30").Value = 475
$ws1.Range("F31").Value = 11010
$ws1.Range("F34").Value = 111
$ws1.Range("F35").Value = 158

# Sheet "全部类型" (sheet4.xml): rows 3-37 in column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3127
$ws4.Range("F5").Value = 1088
$ws4.Range("F7").Value = 30
$ws4.Range("F9").Value = 34
$ws4.Range("F10").Value = 1120
$ws4.Range("F11").Value = 15577
$ws4.Range("F12").Value = 234
$ws4.Range("F13").Value = 168
$ws4.Range("F14").Value = 1022
$ws4.Range("F15").Value = 6142
$ws4.Range("F16").Value = 619
$ws4.Range("F23").Value = 630
$ws4.Range("F24").Value = 9
$ws4.Range("F25").Value = 9
$ws4.Range("F27").Value = 207
$ws4.Range("F28").Value = 860
$ws4.Range("F31").Value = 475
$ws4.Range("F33").Value = 11010
$ws4.Range("F36").Value = 111
$ws4.Range("F37").Value = 158
